$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$sh = $master.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange
Write-Host "Runs: $($tr.Runs().Count)"
for ($i=1; $i -le $tr.Runs().Count; $i++) {
    $r = $tr.Runs($i)
    Write-Host "Run $i text=[$($r.Text)]"
}
